# ============================================================================
# Add a "2022-Q4" quarter sheet (copied from "2022-Q3" as a style template)
# right after "总计" and before "2022-Q3", fill it with the new quarter's
# per-fund holdings, and update the "总计" summary sheet with the new
# 2022-Q4 row (shifting everything else down by one row).
# ============================================================================

function Set-TextValue($ws, $row, $col, $val) {
    # Force the cell to be stored as text even when the string looks like a
    # number (fund codes such as "009073", percentages like "5.62", etc.)
    # so we don't lose leading zeros / trailing zeros to numeric coercion.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-NumValue($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying "2022-Q3" (so it inherits the
#    exact same header/row styling), inserted right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item("2022-Q3 (2)")
$q4Sheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Fill the new sheet's per-fund data (header row already copied as-is).
# ---------------------------------------------------------------------------
$fundRows = @(
    ,@("510810", "汇添富中证上海国企ETF", "101.40", "97.87", "7.02", "7.1183", 4)
    ,@("970007", "华安证券汇赢增利一年持有混合B", "11.05", "22.39", "1.33", "0.1470", 1)
    ,@("970008", "华安证券汇赢增利一年持有混合C", "8.56", "22.39", "1.33", "0.1138", 1)
    ,@("009073", "德邦惠利混合A", "1.12", "53.44", "3.06", "0.0343", 4)
    ,@("001413", "中融鑫起点灵活配置混合A", "0.60", "53.65", "5.47", "0.0328", 2)
    ,@("001739", "中融融安二号灵活配置混合", "0.77", "27.94", "3.81", "0.0293", 3)
    ,@("007924", "方正富邦天鑫灵活配置混合C", "0.45", "45.12", "4.73", "0.0213", 5)
    ,@("009074", "德邦惠利混合C", "0.32", "53.44", "3.06", "0.0098", 4)
    ,@("001414", "中融鑫起点灵活配置混合C", "0.18", "53.65", "5.47", "0.0098", 2)
    ,@("014354", "东方欣冉九个月持有期混合A", "1.13", "23.34", "0.85", "0.0096", 9)
    ,@("014355", "东方欣冉九个月持有期混合C", "1.06", "23.34", "0.85", "0.0090", 9)
    ,@("005373", "中加紫金灵活配置混合A", "0.45", "26.40", "0.62", "0.0028", 6)
    ,@("970006", "华安证券汇赢增利一年持有混合A", "0.18", "22.39", "1.33", "0.0024", 1)
    ,@("005374", "中加紫金灵活配置混合C", "0.07", "26.40", "0.62", "0.0004", 6)
    ,@("007923", "方正富邦天鑫灵活配置混合A", "0.00", "45.12", "4.73", 0, 5)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2   # data starts at row 2
    $row = $fundRows[$i]

    # Column A: the existing 0-based index column. Rows 2-10 already carry
    # the right style (copied from "2022-Q3"); rows 11-16 are new and need
    # the header/index style (same as A10) copied across before writing.
    if ($r -gt 10) {
        $q4Sheet.Range("A10").Copy()
        $q4Sheet.Range("A" + $r).PasteSpecial(-4122)
    }
    Set-NumValue $q4Sheet $r 1 $i

    Set-TextValue $q4Sheet $r 2 $row[0]
    Set-TextValue $q4Sheet $r 3 $row[1]
    Set-TextValue $q4Sheet $r 4 $row[2]
    Set-TextValue $q4Sheet $r 5 $row[3]
    Set-TextValue $q4Sheet $r 6 $row[4]

    # Column G ("持有市值(亿元)") is text except for the very last data row
    # (007923, row 16) which stores a literal numeric 0.
    if ($row[5] -is [string]) {
        Set-TextValue $q4Sheet $r 7 $row[5]
    } else {
        Set-NumValue $q4Sheet $r 7 $row[5]
    }

    Set-NumValue $q4Sheet $r 8 $row[6]
}

$excel.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the
#    top of the data block and shift the existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalRows = @(
    ,@("2022-Q4", 15, 5.62)
    ,@("2022-Q3", 9, 5.61)
    ,@("2022-Q2", 10, 5.44)
    ,@("2022-Q1", 8, 5.21)
    ,@("2021-Q4", 6, 5.34)
    ,@("2021-Q3", 15, 6.21)
    ,@("2021-Q2", 1, 5.35)
    ,@("2021-Q1", 4, 5.5)
    ,@("2020-Q4", 1, 7.12)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    # Row 10 is brand new (the table used to stop at row 9); give its index
    # cell (column A) the same style as the rest of the column before
    # writing the value.
    if ($r -gt 9) {
        $totalSheet.Range("A9").Copy()
        $totalSheet.Range("A" + $r).PasteSpecial(-4122)
    }
    Set-NumValue $totalSheet $r 1 $i

    Set-TextValue $totalSheet $r 2 $row[0]
    Set-NumValue $totalSheet $r 3 $row[1]
    Set-NumValue $totalSheet $r 4 $row[2]
}
$excel.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Keep "2020-Q4" as the selected/active tab, matching the original file.
# ---------------------------------------------------------------------------
$q4OldSheet = $wb.Worksheets.Item("2020-Q4")
$q4OldSheet.Select()
